$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Time" column from G to F (closing the gap left by the
# previously-empty column E/F) and clear out the old column G.
$ws.Range("F1").Value = $ws.Range("G1").Value()
$ws.Range("F2").Value = $ws.Range("G2").Value()
$ws.Range("F3").Value = $ws.Range("G3").Value()
$ws.Range("F4").Value = $ws.Range("G4").Value()
$ws.Range("F5").Value = $ws.Range("G5").Value()

$ws.Range("G1:G5").ClearContents()

# Update the active selection to match the saved state.
$ws.Range("B20").Select()
